$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": remove the two rows for 2023-06-12 and 2023-06-19 ---
# (original rows 3 and 4: 45095.99999999999 / 13 and 45102.99999999999 / 13)
# Deleting these shifts every subsequent row up by two and updates the dimension.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("3:4").Delete()

# --- Sheet "Monthly Trend": requested quantity for the 2023-06 row drops from 34 to 8 ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B3").Value = 8
